$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = 1
$ws.Range("F18").Value = 1

# Row 19
$ws.Range("B19").Value = 1
$ws.Range("F19").Value = 1

# Row 20
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1

# Row 21
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1

# Row 22
$ws.Range("C22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
